# Edit: slide 20 ("Run your Dockerized algorithm") - "Group 8" callout box.
# - Resize/reposition the callout group.
# - Change the rectangle's fill color from amber (FFC000) to red (FF0000).
# - Change the "HINT" callout text to an "ALERT" about running v6 dev
#   create-demo-network on Linux without Docker desktop.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)
$g = $s.Shapes.Item(4)          # Group 8
$rect = $g.GroupItems.Item(1)   # Rectangle 9
$tb   = $g.GroupItems.Item(2)   # TextBox 10

# --- Resize/reposition the whole callout group ------------------------
$g.Left   = 2569945 / 12700
$g.Top    = 5338360 / 12700
$g.Width  = 7016817 / 12700
$g.Height = 1162538 / 12700

# --- Rectangle fill color: FFC000 -> FF0000 ----------------------------
$rect.Fill.ForeColor.RGB = 255

# --- Grow the text box that holds the callout text ---------------------
$tb.Width  = 3805143 / 12700
$tb.Height = 878774 / 12700

# --- Edit the callout text ----------------------------------------------
$tr = $tb.TextFrame.TextRange

# "HINT" -> "ALERT"
$full = $tr.Text
$pos = $full.IndexOf("HINT") + 1
$tr.Characters($pos, 4).Text = "ALERT"

# "Remember to use help(" -> "If you are using Linux without Docker desktop, you need to run:"
$full = $tr.Text
$oldPhrase = "Remember to use help("
$pos = $full.IndexOf($oldPhrase) + 1
$tr.Characters($pos, $oldPhrase.Length).Text = "If you are using Linux without Docker desktop, you need to run:"

# "client.task.create" -> "url"
$full = $tr.Text
$oldPhrase = "client.task.create"
$pos = $full.IndexOf($oldPhrase) + 1
$tr.Characters($pos, $oldPhrase.Length).Text = "url"

# ") to find out which arguments you need to provide" -> " http:172.17.0.1"
$full = $tr.Text
$oldPhrase = ") to find out which arguments you need to provide"
$pos = $full.IndexOf($oldPhrase) + 1
$tr.Characters($pos, $oldPhrase.Length).Text = " http:172.17.0.1"

# Apply the monospace font to the three runs that make up the new command line
$full = $tr.Text
$pos = $full.IndexOf("url") + 1
$tr.Characters($pos, 3).Font.Name = "Aptos Mono"

$full = $tr.Text
$pos = $full.IndexOf(" http:172.17.0.1") + 1
$tr.Characters($pos, 16).Font.Name = "Aptos Mono"

# Insert a new paragraph with the actual v6 dev command, styled as monospace
$tr.InsertAfter([char]13 + "v6 dev create-demo-network -" + [char]8211 + "server-")
$full = $tr.Text
$pos = $full.IndexOf("v6 dev create-demo-network") + 1
$tr.Characters($pos, 37).Font.Name = "Aptos Mono"
